$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change the starting X value from -6.5 to 0 so all subsequent values (built via
# A(n) = A(n-1) + 0.1 formulas) become non-negative.
$ws.Range("A2").Value = 0

# Select cell A2 so it becomes the active cell / selection on the sheet.
$ws.Range("A2").Select() | Out-Null
